$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes old rows 4..48 down to 5..49,
# carrying their existing data/format with them automatically).
$ws.Rows("4:4").Insert()

# New row 4 needs the same formatting as the data rows above/below it
# (bold "A" index column with borders, centered) - copy formats from row 3.
$ws.Range("A3:S3").Copy()
$ws.Range("A4:S4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 index + label: a new "Bruker" scheme inserted right after "Equal Angle".
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Bruker"

# Row 3 ("Equal Angle") values were recalculated.
$ws.Range("C3").Value2 = 0.8952809798270893
$ws.Range("D3").Value2 = 1.138746397694524
$ws.Range("E3").Value2 = 1.013292507204611
$ws.Range("F3").Value2 = 0.8952809798270893
$ws.Range("G3").Value2 = 0.9926729106628243
$ws.Range("H3").Value2 = 1.208206051873199
$ws.Range("I3").Value2 = 0.9762680115273775
$ws.Range("J3").Value2 = 1.138746397694524
$ws.Range("K3").Value2 = 0.8952809798270893
$ws.Range("L3").Value2 = 1.076019452449568
$ws.Range("M3").Value2 = 1.076019452449568
$ws.Range("N3").Value2 = 1.048237271853987
$ws.Range("O3").Value2 = 1.015773294908741
$ws.Range("P3").Value2 = 1.015773294908742
$ws.Range("Q3").Value2 = 0.9856502161383286
$ws.Range("R3").Value2 = 0.9856502161383286
$ws.Range("S3").Value2 = 1.037411143131604

# Row 4 ("Bruker") newly calculated values.
$ws.Range("C4").Value2 = 1.018824535448829
$ws.Range("D4").Value2 = 0.9933260971235289
$ws.Range("E4").Value2 = 0.9937748748443713
$ws.Range("F4").Value2 = 1.018824535448829
$ws.Range("G4").Value2 = 0.9951275235008941
$ws.Range("H4").Value2 = 0.9786944540411573
$ws.Range("I4").Value2 = 0.9973939832842261
$ws.Range("J4").Value2 = 0.9933260971235289
$ws.Range("K4").Value2 = 1.018824535448829
$ws.Range("L4").Value2 = 0.9935504859839501
$ws.Range("M4").Value2 = 0.9935504859839501
$ws.Range("N4").Value2 = 0.9940761651562647
$ws.Range("O4").Value2 = 1.00197516913891
$ws.Range("P4").Value2 = 1.00197516913891
$ws.Range("Q4").Value2 = 1.006187510716389
$ws.Range("R4").Value2 = 1.006187510716389
$ws.Range("S4").Value2 = 0.9961902447071678

